$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

# --- 1. Refresh the "fetched at" timestamp for every existing data row (2-11) ---
$newTimestamp = "2025-12-25 12:37:44"
for ($r = 2; $r -le 11; $r++) {
    $ws.Cells.Item($r, 1).Value = $newTimestamp
}

# --- 2. Insert a new row at position 12, pushing the old row 12 ("限定公開...") down to row 13 ---
$ws.Rows.Item(12).Insert()

# Row 13 (the shifted former row 12) also gets the refreshed timestamp
$ws.Cells.Item(13, 1).Value = $newTimestamp

# --- 3. Populate the newly inserted row 12 with the new job listing ---
$ws.Cells.Item(12, 1).Value = $newTimestamp
$ws.Cells.Item(12, 2).Value = "【急募】Azureサーバー構築の専門家を探しています"
$ws.Cells.Item(12, 3).Value = "システム開発"
$ws.Cells.Item(12, 4).Value = "50,000 円 ~ 100,000 円 / 固定"
$ws.Cells.Item(12, 5).Value = "期限情報なし"
$ws.Cells.Item(12, 6).Value = "https://www.lancers.jp/work/detail/5461140"
$ws.Cells.Item(12, 7).Value = 18

# --- 4. Rebuild the hyperlinks collection ---
# This engine's Hyperlinks.Item(n).Delete() / Range.Hyperlinks (single-cell) don't actually
# detach individual links, and editing .Address on an existing link just appends a stray
# duplicate relationship. The only reliable reset is the whole-sheet delete, so clear
# everything once and re-add every row's link (2-13) fresh, in row order, so the
# relationship ids come out sequential (rId1..rId12) matching F2..F13.
$ws.Range("F2").Hyperlinks.Delete()

$urls = @{
    2  = "https://www.lancers.jp/work/detail/5460562"
    3  = "https://www.lancers.jp/work/detail/5460357"
    4  = "https://www.lancers.jp/work/detail/5460563"
    5  = "https://www.lancers.jp/work/detail/5460750"
    6  = "https://www.lancers.jp/work/detail/5460724"
    7  = "https://www.lancers.jp/work/detail/5460405"
    8  = "https://www.lancers.jp/work/detail/5460928"
    9  = "https://www.lancers.jp/work/detail/5460787"
    10 = "https://www.lancers.jp/work/detail/5016989"
    11 = "https://www.lancers.jp/work/detail/5460484"
    12 = "https://www.lancers.jp/work/detail/5461140"
    13 = "https://www.lancers.jp/work/detail/5450323"
}

for ($r = 2; $r -le 13; $r++) {
    $cell = $ws.Cells.Item($r, 6)
    $ws.Hyperlinks.Add($cell, $urls[$r])
    $cell.Value = $urls[$r]
    $cell.Style = "Hyperlink"
}
